$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 19:22"

# --- Rows whose country stays the same, only the stats refresh ---
$ws.Cells.Item(4,2).Value  = 5268907
$ws.Cells.Item(4,3).Value  = 17461
$ws.Cells.Item(4,4).Value  = 2719664
$ws.Cells.Item(4,5).Value  = 2382430
$ws.Cells.Item(4,7).Value  = 621
$ws.Cells.Item(4,8).Value  = 166813

$ws.Cells.Item(5,2).Value  = 3068138
$ws.Cells.Item(5,3).Value  = 10668
$ws.Cells.Item(5,5).Value  = 802292
$ws.Cells.Item(5,7).Value  = 177
$ws.Cells.Item(5,8).Value  = 102034

$ws.Cells.Item(6,2).Value  = 2325026
$ws.Cells.Item(6,3).Value  = 57873
$ws.Cells.Item(6,4).Value  = 1635812
$ws.Cells.Item(6,5).Value  = 643029
$ws.Cells.Item(6,7).Value  = 832
$ws.Cells.Item(6,8).Value  = 46185

$ws.Cells.Item(13,2).Value = 373692
$ws.Cells.Item(13,3).Value = 3632
$ws.Cells.Item(13,7).Value = 5
$ws.Cells.Item(13,8).Value = 28581

$ws.Cells.Item(21,2).Value = 243180
$ws.Cells.Item(21,3).Value = 1183
$ws.Cells.Item(21,4).Value = 226155
$ws.Cells.Item(21,5).Value = 11152
$ws.Cells.Item(21,7).Value = 15
$ws.Cells.Item(21,8).Value = 5873

$ws.Cells.Item(22,2).Value = 219229
$ws.Cells.Item(22,3).Value = 729
$ws.Cells.Item(22,5).Value = 11062
$ws.Cells.Item(22,7).Value = 2
$ws.Cells.Item(22,8).Value = 9267

$ws.Cells.Item(27,2).Value = 120256
$ws.Cells.Item(27,3).Value = 124
$ws.Cells.Item(27,4).Value = 106524
$ws.Cells.Item(27,5).Value = 4744
$ws.Cells.Item(27,7).Value = 1
$ws.Cells.Item(27,8).Value = 8988

$ws.Cells.Item(54,2).Value = 41404
$ws.Cells.Item(54,3).Value = 192
$ws.Cells.Item(54,4).Value = 39055
$ws.Cells.Item(54,5).Value = 2134

$ws.Cells.Item(59,2).Value = 36204
$ws.Cells.Item(59,3).Value = 492
$ws.Cells.Item(59,4).Value = 25263
$ws.Cells.Item(59,5).Value = 9619
$ws.Cells.Item(59,7).Value = 10
$ws.Cells.Item(59,8).Value = 1322

$ws.Cells.Item(98,2).Value  = 7121
$ws.Cells.Item(98,3).Value  = 309
$ws.Cells.Item(98,5).Value  = 4744
$ws.Cells.Item(98,7).Value  = 7
$ws.Cells.Item(98,8).Value  = 87

$ws.Cells.Item(124,2).Value = 2577
$ws.Cells.Item(124,3).Value = 4
$ws.Cells.Item(124,4).Value = 1973

# --- Rows 66-70: Venezuela/Irlanda/Nepal/Costa Rica/Etiopia re-ranked;
#     each country's row (name + stats) shifts to its new rank position ---
$ws.Cells.Item(66,1).Value = "Irlanda"
$ws.Cells.Item(66,2).Value = 26801
$ws.Cells.Item(66,3).Value = 33
$ws.Cells.Item(66,4).Value = 23364
$ws.Cells.Item(66,5).Value = 1664
$ws.Cells.Item(66,7).Value = 1
$ws.Cells.Item(66,8).Value = 1773

$ws.Cells.Item(67,1).Value = "Venezuela"
$ws.Cells.Item(67,2).Value = 26800
$ws.Cells.Item(67,3).Value = 0
$ws.Cells.Item(67,4).Value = 16930
$ws.Cells.Item(67,5).Value = 9641
$ws.Cells.Item(67,7).Value = 0
$ws.Cells.Item(67,8).Value = 229

$ws.Cells.Item(68,1).Value = "Etiopia"
$ws.Cells.Item(68,2).Value = 24175
$ws.Cells.Item(68,3).Value = 584
$ws.Cells.Item(68,4).Value = 10696
$ws.Cells.Item(68,5).Value = 13039
$ws.Cells.Item(68,7).Value = 20
$ws.Cells.Item(68,8).Value = 440

$ws.Cells.Item(69,1).Value = "Nepal"
$ws.Cells.Item(69,2).Value = 23948
$ws.Cells.Item(69,3).Value = 638
$ws.Cells.Item(69,4).Value = 16664
$ws.Cells.Item(69,5).Value = 7201
$ws.Cells.Item(69,7).Value = 4
$ws.Cells.Item(69,8).Value = 83

$ws.Cells.Item(70,1).Value = "Costa Rica"
$ws.Cells.Item(70,2).Value = 23872
$ws.Cells.Item(70,3).Value = 0
$ws.Cells.Item(70,4).Value = 7823
$ws.Cells.Item(70,5).Value = 15805
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 244

# --- Rows 213-214: Montserrat / Islas Malvinas swap rank ---
$ws.Cells.Item(213,1).Value = "Islas Malvinas"
$ws.Cells.Item(213,4).Value = 13
$ws.Cells.Item(213,8).Value = 0

$ws.Cells.Item(214,1).Value = "Montserrat"
$ws.Cells.Item(214,4).Value = 12
$ws.Cells.Item(214,8).Value = 1
